$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45202 -> 2023-10-03)
# that needs to be bumped by one day (45203 -> 2023-10-04) for every
# data row from row 2 through row 411.
$ws.Range("C2:C411").Value = 45203
